$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("B3").Value = 1.02
$ws.Range("B4").Value = 1.02
$ws.Range("B5").Value = 1.02
$ws.Range("B6").Value = 1.02
$ws.Range("B7").Value = 1.02
$ws.Range("B8").Value = 1.02
$ws.Range("B9").Value = 1.02
$ws.Range("B10").Value = 1.02
$ws.Range("B11").Value = 1.02
$ws.Range("B12").Value = 1.02
$ws.Range("B13").Value = 1.02
$ws.Range("B14").Value = 1.02
$ws.Range("B15").Value = 1.02
$ws.Range("B16").Value = 1.02
$ws.Range("B17").Value = 1.02
$ws.Range("B18").Value = 1.02
$ws.Range("B19").Value = 1.02
$ws.Range("B20").Value = 1.02
$ws.Range("B21").Value = 1.02
$ws.Range("B22").Value = 1.02
$ws.Range("B23").Value = 1.02
$ws.Range("B24").Value = 1.02
$ws.Range("B25").Value = 1.02
$ws.Range("C2").Value = 1.043451830166994
$ws.Range("C3").Value = 1.044288254662686
$ws.Range("C4").Value = 1.0448300831774
$ws.Range("C5").Value = 1.045058011313557
$ws.Range("C6").Value = 1.045096289813848
$ws.Range("C7").Value = 1.044833128202775
$ws.Range("C8").Value = 1.043734377295138
$ws.Range("C9").Value = 1.04180295090964
$ws.Range("C10").Value = 1.04051859951879
$ws.Range("C11").Value = 1.039963256364886
$ws.Range("C12").Value = 1.039757097442194
$ws.Range("C13").Value = 1.03980131376291
$ws.Range("C14").Value = 1.039946212741425
$ws.Range("C15").Value = 1.040035505768002
$ws.Range("C16").Value = 1.040555472588297
$ws.Range("C17").Value = 1.040881846559684
$ws.Range("C18").Value = 1.041072291037848
$ws.Range("C19").Value = 1.041137240534994
$ws.Range("C20").Value = 1.040846821841984
$ws.Range("C21").Value = 1.039903540271035
$ws.Range("C22").Value = 1.039311158261241
$ws.Range("C23").Value = 1.039625124598488
$ws.Range("C24").Value = 1.040862647769469
$ws.Range("C25").Value = 1.042301702379796
$ws.Range("D2").Value = 1.052481821455733
$ws.Range("D3").Value = 1.053155087987962
$ws.Range("D4").Value = 1.053591305905539
$ws.Range("D5").Value = 1.05377482647006
$ws.Range("D6").Value = 1.053805648218919
$ws.Range("D7").Value = 1.053593757589839
$ws.Range("D8").Value = 1.052709236056439
$ws.Range("D9").Value = 1.051155035944626
$ws.Range("D10").Value = 1.050121993663312
$ws.Range("D11").Value = 1.049675430542335
$ws.Range("D12").Value = 1.049509671762139
$ws.Range("D13").Value = 1.049545222366042
$ws.Range("D14").Value = 1.049661726520609
$ws.Range("D15").Value = 1.049733523821396
$ws.Range("D16").Value = 1.050151646564785
$ws.Range("D17").Value = 1.050414126262668
$ws.Range("D18").Value = 1.050567298568877
$ws.Range("D19").Value = 1.050619538581683
$ws.Range("D20").Value = 1.05038595717968
$ws.Range("D21").Value = 1.049627415771503
$ws.Range("D22").Value = 1.049151154887311
$ws.Range("D23").Value = 1.049403566227971
$ws.Range("D24").Value = 1.050398685352477
$ws.Range("D25").Value = 1.051556296854897
$ws.Range("E2").Value = 1.058000762522915
$ws.Range("E3").Value = 1.058863916503
$ws.Range("E4").Value = 1.059423807069321
$ws.Range("E5").Value = 1.059659511647906
$ws.Range("E6").Value = 1.059699106578818
$ws.Range("E7").Value = 1.059426955284603
$ws.Range("E8").Value = 1.058292184163883
$ws.Range("E9").Value = 1.056303167427457
$ws.Range("E10").Value = 1.054984395029488
$ws.Range("E11").Value = 1.054415093354291
$ws.Range("E12").Value = 1.054203892083897
$ws.Range("E13").Value = 1.05424918355284
$ws.Range("E14").Value = 1.054397630020512
$ws.Range("E15").Value = 1.054489127660898
$ws.Range("E16").Value = 1.055022214440566
$ws.Range("E17").Value = 1.055357071730217
$ws.Range("E18").Value = 1.055552555637064
$ws.Range("E19").Value = 1.05561923891673
$ws.Range("E20").Value = 1.055321127391597
$ws.Range("E21").Value = 1.054353908967419
$ws.Range("E22").Value = 1.053747301354502
$ws.Range("E23").Value = 1.054068730691611
$ws.Range("E24").Value = 1.055337368574296
$ws.Range("E25").Value = 1.056816108494728
$ws.Range("F2").Value = 1.065114435861702
$ws.Range("F3").Value = 1.065943137249212
$ws.Range("F4").Value = 1.06648042691526
$ws.Range("F5").Value = 1.066706556390013
$ws.Range("F6").Value = 1.06674453929777
$ws.Range("F7").Value = 1.066483447477075
$ws.Range("F8").Value = 1.06539427781052
$ws.Range("F9").Value = 1.06348325971281
$ws.Range("F10").Value = 1.062214901955727
$ws.Range("F11").Value = 1.061667054600531
$ws.Range("F12").Value = 1.061463765995577
$ws.Range("F13").Value = 1.061507362724271
$ws.Range("F14").Value = 1.061650246481195
$ws.Range("F15").Value = 1.061738309275078
$ws.Range("F16").Value = 1.062251289601644
$ws.Range("F17").Value = 1.062573434097287
$ws.Range("F18").Value = 1.062761466611265
$ws.Range("F19").Value = 1.062825603012455
$ws.Range("F20").Value = 1.062538857469141
$ws.Range("F21").Value = 1.061608165065371
$ws.Range("F22").Value = 1.06102419635679
$ws.Range("F23").Value = 1.061333655327127
$ws.Range("F24").Value = 1.062554480754307
$ws.Range("F25").Value = 1.063976315548619
$ws.Range("I2").Value = 1.044657451029663
$ws.Range("I3").Value = 1.044863234344691
$ws.Range("I4").Value = 1.044995559675705
$ws.Range("I5").Value = 1.045050989760851
$ws.Range("I6").Value = 1.045060284997199
$ws.Range("I7").Value = 1.044996301119953
$ws.Range("I8").Value = 1.044727167933846
$ws.Range("I9").Value = 1.04424659437763
$ws.Range("I10").Value = 1.043922001746932
$ws.Range("I11").Value = 1.043780460067248
$ws.Range("I12").Value = 1.04372773703173
$ws.Range("I13").Value = 1.043739052994461
$ws.Range("I14").Value = 1.043776104982786
$ws.Range("I15").Value = 1.043798914325526
$ws.Range("I16").Value = 1.04393137458593
$ws.Range("I17").Value = 1.044014198617157
$ws.Range("I18").Value = 1.044062412814495
$ws.Range("I19").Value = 1.044078836362487
$ws.Range("I20").Value = 1.044005322271678
$ws.Range("I21").Value = 1.043765198179672
$ws.Range("I22").Value = 1.043613365821046
$ws.Range("I23").Value = 1.04369393600116
$ws.Range("I24").Value = 1.04400933340578
$ws.Range("I25").Value = 1.044371579428462
$ws.Range("J2").Value = 1.04852220459069
$ws.Range("J3").Value = 1.049006042784637
$ws.Range("J4").Value = 1.049319044860052
$ws.Range("J5").Value = 1.049450612067086
$ws.Range("J6").Value = 1.049472701637528
$ws.Range("J7").Value = 1.049320802943204
$ws.Range("J8").Value = 1.048685734460183
$ws.Range("J9").Value = 1.047566154270686
$ws.Range("J10").Value = 1.046819500131128
$ws.Range("J11").Value = 1.046496142934786
$ws.Range("J12").Value = 1.046376027268415
$ws.Range("J13").Value = 1.046401792756245
$ws.Range("J14").Value = 1.046486214259726
$ws.Range("J15").Value = 1.04653822832311
$ws.Range("J16").Value = 1.046840959315753
$ws.Range("J17").Value = 1.047030841697243
$ws.Range("J18").Value = 1.047141591869913
$ws.Range("J19").Value = 1.047179353929007
$ws.Range("J20").Value = 1.047010469623639
$ws.Range("J21").Value = 1.046461354397819
$ws.Range("J22").Value = 1.046116067421395
$ws.Range("J23").Value = 1.04629911361368
$ws.Range("J24").Value = 1.047019674902425
$ws.Range("J25").Value = 1.047855645163283
$ws.Range("K2").Value = 1.055230433028185
$ws.Range("K3").Value = 1.055716912274248
$ws.Range("K4").Value = 1.056031583129814
$ws.Range("K5").Value = 1.056163842202385
$ws.Range("K6").Value = 1.056186047369677
$ws.Range("K7").Value = 1.056033350495877
$ws.Range("K8").Value = 1.055394863722799
$ws.Range("K9").Value = 1.054268955874646
$ws.Range("K10").Value = 1.053517883357141
$ws.Range("K11").Value = 1.053192566515327
$ws.Range("K12").Value = 1.053071716050387
$ws.Range("K13").Value = 1.053097639465902
$ws.Range("K14").Value = 1.053182577241858
$ws.Range("K15").Value = 1.053234908485987
$ws.Range("K16").Value = 1.053539471625062
$ws.Range("K17").Value = 1.053730490882475
$ws.Range("K18").Value = 1.053841899687632
$ws.Range("K19").Value = 1.053879885559677
$ws.Range("K20").Value = 1.053709997296476
$ws.Range("K21").Value = 1.053157565540372
$ws.Range("K22").Value = 1.052810153454197
$ws.Range("K23").Value = 1.052994329960152
$ws.Range("K24").Value = 1.053719257495362
$ws.Range("K25").Value = 1.054560118194346
$ws.Range("L2").Value = 1.060734192414079
$ws.Range("L3").Value = 1.061411168544632
$ws.Range("L4").Value = 1.061849943736719
$ws.Range("L5").Value = 1.062034577113855
$ws.Range("L6").Value = 1.062065587926542
$ws.Range("L7").Value = 1.06185241014199
$ws.Range("L8").Value = 1.060962828301025
$ws.Range("L9").Value = 1.059400905697176
$ws.Range("L10").Value = 1.058363508194452
$ws.Range("L11").Value = 1.057915245225018
$ws.Range("L12").Value = 1.057748882818337
$ws.Range("L13").Value = 1.057784561648099
$ws.Range("L14").Value = 1.057901490741757
$ws.Range("L15").Value = 1.057973553536666
$ws.Range("L16").Value = 1.058393277813321
$ws.Range("L17").Value = 1.058656811785111
$ws.Range("L18").Value = 1.058810616988205
$ws.Range("L19").Value = 1.058863075853474
$ws.Range("L20").Value = 1.058628527732388
$ws.Range("L21").Value = 1.057867054093368
$ws.Range("L22").Value = 1.05738911020196
$ws.Range("L23").Value = 1.057642398569625
$ws.Range("L24").Value = 1.058641307799552
$ws.Range("L25").Value = 1.059804022544505
$ws.Range("M2").Value = 1.067828548415715
$ws.Range("M3").Value = 1.068472552543844
$ws.Range("M4").Value = 1.068889677274683
$ws.Range("M5").Value = 1.069065133306923
$ws.Range("M6").Value = 1.069094598776469
$ws.Range("M7").Value = 1.068892021348438
$ws.Range("M8").Value = 1.068046106565357
$ws.Range("M9").Value = 1.066558715109496
$ws.Range("M10").Value = 1.065569380199614
$ws.Range("M11").Value = 1.065141543333094
$ws.Range("M12").Value = 1.064982710184632
$ws.Range("M13").Value = 1.065016776583571
$ws.Range("M14").Value = 1.065128412402022
$ws.Range("M15").Value = 1.0651972061618
$ws.Range("M16").Value = 1.065597786085719
$ws.Range("M17").Value = 1.065849208138937
$ws.Range("M18").Value = 1.065995911488436
$ws.Range("M19").Value = 1.066045942499894
$ws.Range("M20").Value = 1.065822227425899
$ws.Range("M21").Value = 1.065095536089867
$ws.Range("M22").Value = 1.064639126179341
$ws.Range("M23").Value = 1.064881030739028
$ws.Range("M24").Value = 1.065834418684826
$ws.Range("M25").Value = 1.066942850214103
$ws.Range("N2").Value = 1.050011226388387
$ws.Range("N3").Value = 1.050495751688066
$ws.Range("N4").Value = 1.050809198262331
$ws.Range("N5").Value = 1.050940952309889
$ws.Range("N6").Value = 1.050963073250054
$ws.Range("N7").Value = 1.050810958842162
$ws.Range("N8").Value = 1.050174988489051
$ws.Range("N9").Value = 1.049053818367267
$ws.Range("N10").Value = 1.048306103893189
$ws.Range("N11").Value = 1.047982287492539
$ws.Range("N12").Value = 1.047862001248144
$ws.Range("N13").Value = 1.047887803325923
$ws.Range("N14").Value = 1.047972344717621
$ws.Range("N15").Value = 1.048024432646941
$ws.Range("N16").Value = 1.048327593552318
$ws.Range("N17").Value = 1.048517745588569
$ws.Range("N18").Value = 1.048628653039188
$ws.Range("N19").Value = 1.048666468724737
$ws.Range("N20").Value = 1.048497344584284
$ws.Range("N21").Value = 1.047947449551859
$ws.Range("N22").Value = 1.047601672228324
$ws.Range("N23").Value = 1.047784978367197
$ws.Range("N24").Value = 1.048506562935621
$ws.Range("N25").Value = 1.049343720370138
